$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so numeric-looking strings
# (e.g. "1.001", "28.287.73") are preserved verbatim as text instead of being
# parsed into floating point numbers.
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "28.287.73"
$ws.Range("E2").Value = "  +2.72%  "

$ws.Range("D3").Value = "1.869.27"
$ws.Range("E3").Value = "  +1.30%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").Value = "339.90"
$ws.Range("E5").Value = "  +2.24%  "

$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.12%  "

$ws.Range("D7").Value = "0.4705"
$ws.Range("E7").Value = "  +1.56%  "

$ws.Range("D8").Value = "0.3925"
$ws.Range("E8").Value = "  +2.03%  "

$ws.Range("D9").Value = "47.31"
$ws.Range("E9").Value = "  +2.81%  "

$ws.Range("D10").Value = "0.08003"
$ws.Range("E10").Value = "  +1.21%  "

$ws.Range("D11").Value = "1.006"
$ws.Range("E11").Value = "  +1.25%  "

$ws.Range("D12").Value = "21.87"
$ws.Range("E12").Value = "  +1.95%  "

$ws.Range("D13").Value = "1.879.46"
$ws.Range("E13").Value = "  +1.98%  "

$ws.Range("D14").Value = "6.002"
$ws.Range("E14").Value = "  +1.52%  "

$ws.Range("D15").Value = "7.280"
$ws.Range("E15").Value = "  +2.52%  "

$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Value = "91.20"
$ws.Range("E16").Value = "  +2.67%  "

$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  -0.22%  "

$ws.Range("D18").Value = "0.00001043"
$ws.Range("E18").Value = "  +0.74%  "

$ws.Range("D19").Value = "0.06605"
$ws.Range("E19").Value = "  -0.74%  "

$ws.Range("D20").Value = "17.63"
$ws.Range("E20").Value = "  +3.40%  "

$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  -0.02%  "

$ws.Range("D22").Value = "28.286.04"
$ws.Range("E22").Value = "  +2.71%  "

$ws.Range("D23").Value = "5.448"
$ws.Range("E23").Value = "  +1.28%  "

$ws.Range("D24").Value = "11.06"
$ws.Range("E24").Value = "  +1.46%  "

$ws.Range("D25").Value = "2.289"
$ws.Range("E25").Value = "  -0.61%  "

$ws.Range("D26").Value = "2.093.82"
$ws.Range("E26").Value = "  +1.49%  "

$ws.Range("D27").Value = "159.92"
$ws.Range("E27").Value = "  +1.23%  "

$ws.Range("D28").Value = "19.87"
$ws.Range("E28").Value = "  +2.03%  "

$ws.Range("D29").Value = "2.147"
$ws.Range("E29").Value = "  +2.26%  "

$ws.Range("D30").Value = "5.507"
$ws.Range("E30").Value = "  +2.09%  "

$ws.Range("D31").Value = "120.25"
$ws.Range("E31").Value = "  +0.52%  "

$ws.Range("D32").Value = "0.9780"
$ws.Range("E32").Value = "  +0.35%  "

$ws.Range("D33").Value = "0.09519"
$ws.Range("E33").Value = "  +1.21%  "

$ws.Range("D34").Value = "3.593"
$ws.Range("E34").Value = "  +0.53%  "

$ws.Range("D35").Value = "1.376"
$ws.Range("E35").Value = "  +2.38%  "

$ws.Range("D36").Value = "5.353"
$ws.Range("E36").Value = "  +1.44%  "

$ws.Range("D37").Value = "0.02269"
$ws.Range("E37").Value = "  +2.20%  "

$ws.Range("D38").Value = "0.06096"
$ws.Range("E38").Value = "  +1.41%  "

$ws.Range("D39").Value = "8.426"
$ws.Range("E39").Value = "  +1.79%  "

$ws.Range("D40").Value = "1.179"
$ws.Range("E40").Value = "  +0.06%  "

$ws.Range("D41").Value = "0.5976"
$ws.Range("E41").Value = "  +1.49%  "

$ws.Range("E42").Value = "  +0.01%  "

$ws.Range("D43").Value = "0.1882"
$ws.Range("E43").Value = "  +1.23%  "

$ws.Range("D44").Value = "10.39"
$ws.Range("E44").Value = "  +1.16%  "

$ws.Range("D45").Value = "1.285"
$ws.Range("E45").Value = "  +2.57%  "

$ws.Range("D46").Value = "0.5616"
$ws.Range("E46").Value = "  +0.69%  "

$ws.Range("D47").Value = "12.15"
$ws.Range("E47").Value = "  -0.39%  "

$ws.Range("D48").Value = "1.968"
$ws.Range("E48").Value = "  +3.83%  "

$ws.Range("D49").Value = "0.06906"
$ws.Range("E49").Value = "  +3.32%  "

$ws.Range("D50").Value = "111.32"
$ws.Range("E50").Value = "  +0.67%  "

$ws.Range("D51").Value = "2.025"
$ws.Range("E51").Value = "  +13.65%  "

# Restore the default cell style so no stray number-format styling is
# introduced (matches the original formatting of these cells).
$priceRange.Style = "Normal"
